# Weekly update: insert a new record as row 3 (pushing all existing data
# rows down by one) and let the last existing row (old 108) shift to the
# new row 109 — matching the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 3..108 down to 4..109 (this also bumps the sheet dimension
# from R108 to R109 and carries the D-column date style along).
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with this week's record.
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 44922
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 100112031
$ws.Range("G3").Value = "Poroto verde"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 32000
$ws.Range("L3").Value = 32000
$ws.Range("M3").Value = 32000
$ws.Range("N3").Value = "$/saco 25 kilos"
$ws.Range("O3").Value = "Región del Maule"
$ws.Range("P3").Value = 1280
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
